# Update weekly Fruta/Hortaliza (Papaya) price records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44400
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 1500
$ws.Range("O2").Value = 1500
$ws.Range("P2").Value = 1500
$ws.Range("S2").Value = 1500

# Row 3
$ws.Range("D3").Value = 44292
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 14000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1400
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44195
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1500
$ws.Range("T4").Value = 10

# Row 5
$ws.Range("D5").Value = 44391
$ws.Range("M5").Value = 15
$ws.Range("N5").Value = 1500
$ws.Range("O5").Value = 1500
$ws.Range("P5").Value = 1500
$ws.Range("Q5").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("T5").Value = 1

# Row 6
$ws.Range("D6").Value = 44391
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 1000
$ws.Range("O6").Value = 1000
$ws.Range("P6").Value = 1000
$ws.Range("Q6").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 1

# Row 10
$ws.Range("D10").Value = 44343
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 1700
$ws.Range("O10").Value = 1700
$ws.Range("P10").Value = 1700
$ws.Range("S10").Value = 1700

# Row 11
$ws.Range("D11").Value = 44309
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 1600
$ws.Range("O11").Value = 1600
$ws.Range("P11").Value = 1600
$ws.Range("S11").Value = 1600
